# Restore cell C10 on the "Rules" sheet to the numeric value 1
# (was 18, author's commit message indicates a restore of a previous revision).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1

